# RPA datasets push 2024-07-19
# Insert the new "티디에스팜" demand-forecast row at row 5 (in date order),
# shifting the existing rows 5-14 down by one, then drop the now-duplicated
# old row (which lands on row 15) so the remainder of the table (old rows
# 15-21 / old "아이빔테크놀로지" onward) lines back up unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 5, pushing rows 5..21 down to 6..22.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new offering's data.
$ws.Range("A5").Value = "티디에스팜"
$ws.Range("B5").Value = "2024.07.31~08.06"
$ws.Range("C5").Value = "9,500~10,700"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 9500
$ws.Range("F5").Value = "한국투자증권"

# Row 15 now holds the stale duplicate (old row 14's "티디에스팜" data);
# remove it so everything below collapses back up by one row.
$ws.Rows("15:15").Delete()
